$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "H", "J", "L", "N", "O")

$data = @(
    @(19.22748317849181, 8.355915025247411, 8.103667696083248, 12.63965890866636, 35.63091631070544, 7.344005520526261, 9.859911517634377, 11.19252879243041, 18.51592985720283, 27.36555060436836),
    @(18.8310529769231, 8.125909172112635, 8.10361583943963, 12.66858514225142, 35.68512573369716, 7.344005520526261, 9.882424730450563, 11.17846103245517, 18.57153812877309, 27.42774904111684),
    @(18.58641455618903, 7.980213968520701, 8.10441478491089, 12.68773708461488, 35.72695103253579, 7.344005520526261, 9.897002703885537, 11.17119862931815, 18.60755564933989, 27.47217397885904),
    @(18.48655056470074, 7.919775578427736, 8.104950050786941, 12.69589192634582, 35.74613865516616, 7.344005520526261, 9.903133647265431, 11.16858732881962, 18.62270535606318, 27.49184163534586),
    @(18.46996192502061, 7.909677198482005, 8.105051617930927, 12.69726720059509, 35.74945407147543, 7.344005520526261, 9.904163194405664, 11.16817481990404, 18.62524950981157, 27.49520179084912),
    @(18.58506826740705, 7.979403113445644, 8.104421153664235, 12.68784564481635, 35.72720113068282, 7.344005520526261, 9.897084616736054, 11.17116199960381, 18.60775804991535, 27.47243289575747),
    @(19.09113009239249, 8.277570153677489, 8.103477866897506, 12.64934422573422, 35.64783288696815, 7.344005520526261, 9.867517727160255, 11.18739426166459, 18.53471533193919, 27.38570047220432),
    @(20.06796052008097, 8.82449501727042, 8.108184646511269, 12.58486036538966, 35.56011127419439, 7.344005520526261, 9.815502693231327, 11.23002794669172, 18.40629728028926, 27.26524078556178),
    @(20.76874861954599, 9.200517073293325, 8.11559079084652, 12.54417237425893, 35.53723489104181, 7.344005520526261, 9.78089146065328, 11.2677806887405, 18.32091351637763, 27.2071786534323),
    @(21.08244224020467, 9.365476359181145, 8.119805954673129, 12.5271086934959, 35.53587382266505, 7.344005520526261, 9.765921632414369, 11.2863150046207, 18.28400255283729, 27.18740703136162),
    @(21.2003854121645, 9.427029319143752, 8.12152272676904, 12.52085452089929, 35.53665913752381, 7.344005520526261, 9.76036387477008, 11.29352569144129, 18.27030190274303, 27.18087698082807),
    @(21.1750235731461, 9.413813967951681, 8.12114764409943, 12.52219224761547, 35.53643216589519, 7.344005520526261, 9.761555908091642, 11.29196425253497, 18.2732402889436, 27.18224075180026),
    @(21.0921630663607, 9.370558865849647, 8.119944784308391, 12.52659000310067, 35.53591236591197, 7.344005520526261, 9.765462170486138, 11.28690440183679, 18.28286985096265, 27.18685060963305),
    @(21.04129524655145, 9.343943861340009, 8.119223668149161, 12.52931076294367, 35.53576334982318, 7.344005520526261, 9.767869310022983, 11.28383001386796, 18.28880424857854, 27.18979896795032),
    @(20.74813561609268, 9.189610690609646, 8.115332253692866, 12.54531657982824, 35.53750587875537, 7.344005520526261, 9.781885330274228, 11.26659650861685, 18.32336450641272, 27.2086045969476),
    @(20.56690341627054, 9.093343592225379, 8.113160941916453, 12.55550557964415, 35.54089191781542, 7.344005520526261, 9.790681876461962, 11.25637031919738, 18.34505993020115, 27.22184378925772),
    @(20.46218984696122, 9.037401574655437, 8.111991765148598, 12.561502102389, 35.54369094256854, 7.344005520526261, 9.795814390079736, 11.2506167520367, 18.35772032661955, 27.23008362216783),
    @(20.42665784089751, 9.018363566856253, 8.111609623433438, 12.56355580553148, 35.54478486958683, 7.344005520526261, 9.797564718021658, 11.2486908263188, 18.36203816435019, 27.23298076283866),
    @(20.58624574825544, 9.103650841750349, 8.113383840351093, 12.55440686230273, 35.54044334193526, 7.344005520526261, 9.789737919968458, 11.25744566411095, 18.34273160999226, 27.22036975327363),
    @(21.11652502214944, 9.383289007766466, 8.120294829737613, 12.52529264803816, 35.53602974718223, 7.344005520526261, 9.764311797834141, 11.28838541646442, 18.28003391166698, 27.1854705948886),
    @(21.45811613888697, 9.560709332482991, 8.12551386089641, 12.50747395768527, 35.54072615581297, 7.344005520526261, 9.748341093843967, 11.30972456494274, 18.24066994858188, 27.16824093887675),
    @(21.27629326606183, 9.46651659771794, 8.122664476235588, 12.51687362026422, 35.53752617812466, 7.344005520526261, 9.756805930751391, 11.29823433327508, 18.26153196710732, 27.17692569308176),
    @(20.57750269335235, 9.09899279026607, 8.113282821389946, 12.55490316025446, 35.54064348806568, 7.344005520526261, 9.790164448227157, 11.25695910916946, 18.34378365970834, 27.2210342075002),
    @(19.80615498248554, 8.680885141230391, 8.106214270931869, 12.60112839112652, 35.57655051307722, 7.344005520526261, 9.828938867120472, 11.21735331956632, 18.43945853815866, 27.2924949401667)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellRef = "{0}{1}" -f $cols[$j], $rowNum
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
